$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 124.85714
$ws.Range("I11").Value = 124.85714
$ws.Range("K11").Value = 124.85714
$ws.Range("M11").Value = 15.14286
$ws.Range("H17").Value = 2519.7646
$ws.Range("I17").Value = 3109.5
$ws.Range("K17").Value = 9328.5
$ws.Range("M17").Value = -9160.5
$ws.Range("H28").Value = 453.75
$ws.Range("I28").Value = 453.75
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 453.75
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 31.25
$ws.Range("N28").ClearContents()
$ws.Range("H40").Value = 2307.6155
$ws.Range("I40").Value = 1999.875
$ws.Range("K40").Value = 1999.875
$ws.Range("M40").Value = -1824.875
$ws.Range("H42").Value = 120
$ws.Range("I42").Value = 102.2
$ws.Range("J42").Value = 209
$ws.Range("K42").Value = 306.6
$ws.Range("L42").Value = 627
$ws.Range("M42").Value = -76.60000000000002
$ws.Range("N42").Value = -1087
$ws.Range("H62").Value = 5984.143
$ws.Range("I62").Value = 3968.2856
$ws.Range("K62").Value = 3968.2856
$ws.Range("M62").Value = -3344.2856
$ws.Range("H65").Value = 5984.143
$ws.Range("I65").Value = 3968.2856
$ws.Range("K65").Value = 19841.428
$ws.Range("M65").Value = -16721.428
$ws.Range("H76").Value = 5814.1304
$ws.Range("I76").Value = 4671.5386
$ws.Range("K76").Value = 4671.5386
$ws.Range("M76").Value = -4356.5386
$ws.Range("H79").Value = 5814.1304
$ws.Range("I79").Value = 4671.5386
$ws.Range("K79").Value = 4671.5386
$ws.Range("M79").Value = -3579.5386
$ws.Range("H86").Value = 5296.5
$ws.Range("J86").Value = 5638.4287
$ws.Range("L86").Value = 5638.4287
$ws.Range("N86").Value = -7884.4287
$ws.Range("H89").Value = 5296.5
$ws.Range("J89").Value = 5638.4287
$ws.Range("L89").Value = 28192.1435
$ws.Range("N89").Value = -39424.14350000001
$ws.Range("H98").Value = 1101.2858
$ws.Range("I98").Value = 1101.2858
$ws.Range("K98").Value = 1101.2858
$ws.Range("M98").Value = 396.7141999999999
$ws.Range("H106").Value = 32180
$ws.Range("I106").Value = 36993.832
$ws.Range("J106").Value = 3297
$ws.Range("K106").Value = 36993.832
$ws.Range("L106").Value = 3297
$ws.Range("M106").Value = -36362.832
$ws.Range("N106").Value = -4559
$ws.Range("H122").Value = 1101.2858
$ws.Range("I122").Value = 1101.2858
$ws.Range("K122").Value = 3303.8574
$ws.Range("M122").Value = -853.8574000000003
$ws.Range("H132").Value = 1315.5428
$ws.Range("I132").Value = 1254.3636
$ws.Range("K132").Value = 3763.0908
$ws.Range("M132").Value = -1233.0908
$ws.Range("H138").Value = 3404.9253
$ws.Range("J138").Value = 3619.9673
$ws.Range("L138").Value = 10859.9019
$ws.Range("N138").Value = -21139.9019

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2499
$ws.Range("I45").Value = 2499
$ws.Range("K45").Value = 2499
$ws.Range("M45").Value = -2122
$ws.Range("H74").Value = 1607.5714
$ws.Range("I74").Value = 841.2105
$ws.Range("K74").Value = 841.2105
$ws.Range("M74").Value = 32.78949999999998
$ws.Range("H77").Value = 1607.5714
$ws.Range("I77").Value = 841.2105
$ws.Range("K77").Value = 4206.0525
$ws.Range("M77").Value = 161.9475000000002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1731.5385
$ws.Range("J20").Value = 15000
$ws.Range("L20").Value = 15000
$ws.Range("N20").Value = -15494
$ws.Range("H80").Value = 347.16666
$ws.Range("I80").Value = 165
$ws.Range("J80").Value = 529.3333
$ws.Range("K80").Value = 165
$ws.Range("L80").Value = 529.3333
$ws.Range("M80").Value = 833
$ws.Range("N80").Value = -2525.3333
$ws.Range("H83").Value = 347.16666
$ws.Range("I83").Value = 165
$ws.Range("J83").Value = 529.3333
$ws.Range("K83").Value = 825
$ws.Range("L83").Value = 2646.6665
$ws.Range("M83").Value = 4167
$ws.Range("N83").Value = -12630.6665
$ws.Range("H99").Value = 3190.3076
$ws.Range("I99").Value = 3019.6956
$ws.Range("K99").Value = 3019.6956
$ws.Range("M99").Value = -1521.6956

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 296.8
$ws.Range("I7").Value = 296.8
$ws.Range("K7").Value = 296.8
$ws.Range("M7").Value = -183.8
$ws.Range("H16").Value = 1783.6
$ws.Range("I16").Value = 1479.75
$ws.Range("K16").Value = 1479.75
$ws.Range("M16").Value = -1192.75
$ws.Range("H86").Value = 11951.857
$ws.Range("I86").Value = 9422
$ws.Range("K86").Value = 9422
$ws.Range("M86").Value = -8299
$ws.Range("H89").Value = 11951.857
$ws.Range("I89").Value = 9422
$ws.Range("K89").Value = 47110
$ws.Range("M89").Value = -41494
$ws.Range("H113").Value = 1783.6
$ws.Range("I113").Value = 1479.75
$ws.Range("K113").Value = 1479.75
$ws.Range("M113").Value = 690.25
$ws.Range("H122").Value = 3242.5386
$ws.Range("I122").Value = 3488.72
$ws.Range("J122").Value = 2802.9285
$ws.Range("K122").Value = 10466.16
$ws.Range("L122").Value = 8408.7855
$ws.Range("M122").Value = -8016.16
$ws.Range("N122").Value = -13308.7855
$ws.Range("H134").Value = 1953.9286
$ws.Range("I134").Value = 1014.1905
$ws.Range("K134").Value = 3042.5715
$ws.Range("M134").Value = -507.5715
$ws.Range("H141").Value = 90642.086
$ws.Range("J141").Value = 90642.086
$ws.Range("L141").Value = 90642.086
$ws.Range("N141").Value = -101002.086

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 650.5
$ws.Range("I14").Value = 650.5
$ws.Range("K14").Value = 1951.5
$ws.Range("M14").Value = -1778.5
$ws.Range("H137").Value = 9917.799999999999
$ws.Range("I137").Value = 9749.5
$ws.Range("J137").Value = 9959.875
$ws.Range("K137").Value = 29248.5
$ws.Range("L137").Value = 29879.625
$ws.Range("M137").Value = -24148.5
$ws.Range("N137").Value = -40079.625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H113").Value = 4254.3076
$ws.Range("I113").Value = 1623.75
$ws.Range("K113").Value = 1623.75
$ws.Range("M113").Value = 546.25
$ws.Range("H122").Value = 5141.4546
$ws.Range("I122").Value = 4642.1113
$ws.Range("J122").Value = 7388.5
$ws.Range("K122").Value = 13926.3339
$ws.Range("L122").Value = 22165.5
$ws.Range("M122").Value = -11476.3339
$ws.Range("N122").Value = -27065.5
$ws.Range("H132").Value = 4326.1113
$ws.Range("I132").Value = 2652
$ws.Range("K132").Value = 7956
$ws.Range("M132").Value = -5426

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2394.6667
$ws.Range("I22").Value = 1873.6
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 1873.6
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -1578.6
$ws.Range("N22").Value = -5590
$ws.Range("H27").Value = 2394.6667
$ws.Range("I27").Value = 1873.6
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 1873.6
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = -1766.6
$ws.Range("N27").Value = -5214
$ws.Range("H122").Value = 20001
$ws.Range("I122").Value = 17501.5
$ws.Range("J122").Value = 25000
$ws.Range("K122").Value = 52504.5
$ws.Range("L122").Value = 75000
$ws.Range("M122").Value = -50054.5
$ws.Range("N122").Value = -79900
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1079
$ws.Range("I113").Value = 983.6
$ws.Range("J113").Value = 1397
$ws.Range("K113").Value = 2950.8
$ws.Range("L113").Value = 4191
$ws.Range("M113").Value = -780.8000000000002
$ws.Range("N113").Value = -8531
$ws.Range("H132").Value = 1408.9
$ws.Range("I132").Value = 1246.2106
$ws.Range("K132").Value = 3738.6318
$ws.Range("M132").Value = -1208.6318
$ws.Range("H136").Value = 4424.231
$ws.Range("I136").Value = 1283.7142
$ws.Range("J136").Value = 8088.1665
$ws.Range("K136").Value = 3851.1426
$ws.Range("L136").Value = 24264.4995
$ws.Range("M136").Value = -1301.1426
$ws.Range("N136").Value = -29364.4995
